# Scheduled data refresh: updates market-price-derived figures (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Values below come from an external pricing source; cells whose computed profit
# collapses to zero/blank are cleared instead of written as 0.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1492.8462
$ws.Range("I38").Value = 40.7
$ws.Range("K38").Value = 122.1
$ws.Range("M38").Value = 249.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 13845.571
$ws.Range("I69").Value = 10985.363
$ws.Range("K69").Value = 32956.089
$ws.Range("M69").Value = -32082.089

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 13845.571
$ws.Range("I72").Value = 10985.363
$ws.Range("K72").Value = 98868.26699999999
$ws.Range("M72").Value = -94500.26699999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8289.815000000001
$ws.Range("I132").Value = 4481.394
$ws.Range("J132").Value = 33425.4
$ws.Range("K132").Value = 13444.182
$ws.Range("L132").Value = 100276.2
$ws.Range("M132").Value = -10914.182
$ws.Range("N132").Value = -105336.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2747.8572
$ws.Range("I138").Value = 1780.1875
$ws.Range("K138").Value = 5340.5625
$ws.Range("M138").Value = -200.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1817.2174
$ws.Range("I2").Value = 1737.3158
$ws.Range("K2").Value = 1737.3158
$ws.Range("M2").Value = -1624.3158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10837.404
$ws.Range("I32").Value = 7846.8613
$ws.Range("K32").Value = 7846.8613
$ws.Range("M32").Value = -7559.8613

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8357.261
$ws.Range("I45").Value = 9136.368
$ws.Range("J45").Value = 4656.5
$ws.Range("K45").Value = 9136.368
$ws.Range("L45").Value = 4656.5
$ws.Range("M45").Value = -8759.368
$ws.Range("N45").Value = -5410.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4930.7896
$ws.Range("I61").Value = 4772.2583
$ws.Range("J61").Value = 5632.857
$ws.Range("K61").Value = 4772.2583
$ws.Range("L61").Value = 5632.857
$ws.Range("M61").Value = -4560.2583
$ws.Range("N61").Value = -6056.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 660.3077
$ws.Range("I74").Value = 615.5833
$ws.Range("J74").Value = 1197
$ws.Range("K74").Value = 615.5833
$ws.Range("L74").Value = 1197
$ws.Range("M74").Value = 258.4167
$ws.Range("N74").Value = -2945

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 660.3077
$ws.Range("I77").Value = 615.5833
$ws.Range("J77").Value = 1197
$ws.Range("K77").Value = 3077.9165
$ws.Range("L77").Value = 5985
$ws.Range("M77").Value = 1290.0835
$ws.Range("N77").Value = -14721

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1817.2174
$ws.Range("I116").Value = 1737.3158
$ws.Range("K116").Value = 1737.3158
$ws.Range("M116").Value = 556.6841999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 154999.5
$ws.Range("J133").Value = 154999.5
$ws.Range("L133").Value = 154999.5
$ws.Range("N133").Value = -160059.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4930.7896
$ws.Range("I136").Value = 4772.2583
$ws.Range("J136").Value = 5632.857
$ws.Range("K136").Value = 14316.7749
$ws.Range("L136").Value = 16898.571
$ws.Range("M136").Value = -11766.7749
$ws.Range("N136").Value = -21998.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1817.2174
$ws.Range("I3").Value = 1737.3158
$ws.Range("K3").Value = 1737.3158
$ws.Range("M3").Value = -1623.3158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1371.9
$ws.Range("J94").Value = 1806.3334
$ws.Range("L94").Value = 1806.3334
$ws.Range("N94").Value = -2708.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2961.5
$ws.Range("I107").Value = 2786
$ws.Range("K107").Value = 2786
$ws.Range("M107").Value = -866

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13581.28
$ws.Range("I99").Value = 9466.083000000001
$ws.Range("K99").Value = 9466.083000000001
$ws.Range("M99").Value = -7968.083000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5223.6284
$ws.Range("I122").Value = 7860.706
$ws.Range("J122").Value = 2733.0557
$ws.Range("K122").Value = 23582.118
$ws.Range("L122").Value = 8199.167099999999
$ws.Range("M122").Value = -21132.118
$ws.Range("N122").Value = -13099.1671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13581.28
$ws.Range("I126").Value = 9466.083000000001
$ws.Range("K126").Value = 28398.249
$ws.Range("M126").Value = -25928.249

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5192.6514
$ws.Range("I132").Value = 3152.8125
$ws.Range("J132").Value = 11126.728
$ws.Range("K132").Value = 9458.4375
$ws.Range("L132").Value = 33380.18399999999
$ws.Range("M132").Value = -6928.4375
$ws.Range("N132").Value = -38440.18399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1983.5686
$ws.Range("I134").Value = 1983.5686
$ws.Range("K134").Value = 5950.7058
$ws.Range("M134").Value = -3415.7058

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 91820820
$ws.Range("J9").Value = 111114220
$ws.Range("L9").Value = 333342660
$ws.Range("N9").Value = -333343108

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 5332
$ws.Range("J57").Value = 8000
$ws.Range("L57").Value = 24000
$ws.Range("N57").Value = -25118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 8499.333000000001
$ws.Range("J88").Value = 8499.333000000001
$ws.Range("L88").Value = 25497.999
$ws.Range("N88").Value = -26353.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 8499.333000000001
$ws.Range("J91").Value = 8499.333000000001
$ws.Range("L91").Value = 25497.999
$ws.Range("N91").Value = -28461.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 647.7778
$ws.Range("I113").Value = 647.7778
$ws.Range("K113").Value = 1943.3334
$ws.Range("M113").Value = 226.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2427.25
$ws.Range("J132").Value = 2782.9412
$ws.Range("L132").Value = 25046.4708
$ws.Range("N132").Value = -30106.4708

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2054.6924
$ws.Range("J137").Value = 1450
$ws.Range("L137").Value = 4350
$ws.Range("N137").Value = -14550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10955444
$ws.Range("I80").Value = 19168734
$ws.Range("K80").Value = 19168734
$ws.Range("M80").Value = -19167736

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 10955444
$ws.Range("I83").Value = 19168734
$ws.Range("K83").Value = 95843670
$ws.Range("M83").Value = -95838678

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2838.7368
$ws.Range("I113").Value = 2415.3076
$ws.Range("K113").Value = 2415.3076
$ws.Range("M113").Value = -245.3076000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6863.077
$ws.Range("I126").Value = 5313.143
$ws.Range("J126").Value = 8671.333000000001
$ws.Range("K126").Value = 15939.429
$ws.Range("L126").Value = 26013.999
$ws.Range("M126").Value = -13469.429
$ws.Range("N126").Value = -30953.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7712.846
$ws.Range("I40").Value = 7661.727
$ws.Range("K40").Value = 7661.727
$ws.Range("M40").Value = -7525.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3253
$ws.Range("I46").Value = 1079
$ws.Range("J46").Value = 4611.75
$ws.Range("K46").Value = 1079
$ws.Range("L46").Value = 4611.75
$ws.Range("M46").Value = -891
$ws.Range("N46").Value = -4987.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2582
$ws.Range("I93").Value = 2934.3333
$ws.Range("J93").Value = 1525
$ws.Range("K93").Value = 2934.3333
$ws.Range("L93").Value = 1525
$ws.Range("M93").Value = -1686.3333
$ws.Range("N93").Value = -4021

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 12099
$ws.Range("J101").Value = 12099
$ws.Range("L101").Value = 12099
$ws.Range("N101").Value = -18589

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 17422.555
$ws.Range("I136").Value = 35834.668
$ws.Range("J136").Value = 8216.5
$ws.Range("K136").Value = 107504.004
$ws.Range("L136").Value = 24649.5
$ws.Range("M136").Value = -104954.004
$ws.Range("N136").Value = -29749.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 43166
$ws.Range("I37").Value = 49749
$ws.Range("K37").Value = 49749
$ws.Range("M37").Value = -49546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 33332
$ws.Range("J103").Value = 33332
$ws.Range("L103").Value = 33332
$ws.Range("N103").Value = -35676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4538.979
$ws.Range("I122").Value = 2252.225
$ws.Range("K122").Value = 6756.674999999999
$ws.Range("M122").Value = -4306.674999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2544.8723
$ws.Range("I132").Value = 2758.9062
$ws.Range("J132").Value = 2088.2666
$ws.Range("K132").Value = 8276.7186
$ws.Range("L132").Value = 6264.7998
$ws.Range("M132").Value = -5746.7186
$ws.Range("N132").Value = -11324.7998
